$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Alimentos")

# --- Fix H/I columns for existing rows (rows 2-16) ---
$ws.Range("H2").Value = 10
$ws.Range("I2").Value = 0
$ws.Range("H3").Value = 26
$ws.Range("I3").Value = 0
$ws.Range("H4").Value = 26
$ws.Range("I4").Value = 0
$ws.Range("H5").Value = 10
$ws.Range("I5").Value = 0
$ws.Range("H6").Value = 5
$ws.Range("I6").Value = 1
$ws.Range("H7").Value = 10
$ws.Range("I7").Value = 1
$ws.Range("H8").Value = 16
$ws.Range("I8").Value = 0
$ws.Range("H9").Value = 16
$ws.Range("I9").Value = 0
$ws.Range("H10").Value = 16
$ws.Range("I10").Value = 0
$ws.Range("H11").Value = 16
$ws.Range("I11").Value = 1
$ws.Range("H12").Value = 16
$ws.Range("I12").Value = 0
$ws.Range("H13").Value = 16
$ws.Range("I13").Value = 0
$ws.Range("H14").Value = 16
$ws.Range("I14").Value = 0
$ws.Range("H15").Value = 26
$ws.Range("I15").Value = 0
$ws.Range("H16").Value = 7
$ws.Range("I16").Value = 1

# --- Add new rows 17-29 ---
$ws.Range("A17").Value = '100 gr arroz + 200 gr filetes pollo empanados + yogurt de chocolate'
$ws.Range("B17").Value = 1018
$ws.Range("C17").Value = 37.54
$ws.Range("D17").Value = 8.51
$ws.Range("E17").Value = 121.38
$ws.Range("F17").Value = 8.640000000000001
$ws.Range("G17").Value = 46.44
$ws.Range("H17").Value = 5
$ws.Range("I17").Value = 0
$ws.Range("J17").Value = 3

$ws.Range("A18").Value = 'pizza precocinada mediana '
$ws.Range("B18").Value = 844
$ws.Range("C18").Value = 39.2
$ws.Range("D18").Value = 16.48
$ws.Range("E18").Value = 88.8
$ws.Range("F18").Value = 0
$ws.Range("G18").Value = 32.8
$ws.Range("H18").Value = 5
$ws.Range("I18").Value = 0
$ws.Range("J18").Value = 4

$ws.Range("A19").Value = '100 GR ARROZ + 200 gr pollo frito + manzana'
$ws.Range("B19").Value = 689
$ws.Range("C19").Value = 15.51
$ws.Range("D19").Value = 2.88
$ws.Range("E19").Value = 110.74
$ws.Range("F19").Value = 24.2
$ws.Range("G19").Value = 24.38
$ws.Range("H19").Value = 5
$ws.Range("I19").Value = 0
$ws.Range("J19").Value = 1

$ws.Range("A20").Value = '150 gr lomo de ternera + 5 croquetas + racion patatas fritas'
$ws.Range("B20").Value = 821.2
$ws.Range("C20").Value = 31.066
$ws.Range("D20").Value = 4.772
$ws.Range("E20").Value = 82.419
$ws.Range("F20").Value = 8.085000000000001
$ws.Range("G20").Value = 51.27399999999999
$ws.Range("H20").Value = 5
$ws.Range("I20").Value = 0
$ws.Range("J20").Value = 4

$ws.Range("A21").Value = '150 gr de pasta + 100 gr lomo de ternera + platano'
$ws.Range("B21").Value = 781.5
$ws.Range("C21").Value = 5.799999999999999
$ws.Range("D21").Value = 1.465
$ws.Range("E21").Value = 136.35
$ws.Range("F21").Value = 0
$ws.Range("G21").Value = 40.75
$ws.Range("H21").Value = 5
$ws.Range("I21").Value = 1
$ws.Range("J21").Value = 2

$ws.Range("A22").Value = 'manzana'
$ws.Range("B22").Value = 100
$ws.Range("C22").Value = 0
$ws.Range("D22").Value = 0
$ws.Range("E22").Value = 24
$ws.Range("F22").Value = 24
$ws.Range("G22").Value = 0.3
$ws.Range("H22").Value = 26
$ws.Range("I22").Value = 0
$ws.Range("J22").Value = 1

$ws.Range("A23").Value = 'platano'
$ws.Range("B23").Value = 140
$ws.Range("C23").Value = 0.4
$ws.Range("D23").Value = 0.18
$ws.Range("E23").Value = 30
$ws.Range("G23").Value = 1.8
$ws.Range("H23").Value = 26
$ws.Range("I23").Value = 0
$ws.Range("J23").Value = 1

$ws.Range("A24").Value = 'tostada de atun (50 gr)'
$ws.Range("B24").Value = 242.5
$ws.Range("C24").Value = 1.54
$ws.Range("D24").Value = 0.095
$ws.Range("E24").Value = 48
$ws.Range("F24").Value = 0
$ws.Range("G24").Value = 8.470000000000001
$ws.Range("H24").Value = 10
$ws.Range("I24").Value = 0
$ws.Range("J24").Value = 1

$ws.Range("A25").Value = 'yogurt griego'
$ws.Range("B25").Value = 139
$ws.Range("C25").Value = 10.2
$ws.Range("D25").Value = 6.8
$ws.Range("E25").Value = 5.4
$ws.Range("F25").Value = 5.3
$ws.Range("G25").Value = 6.4
$ws.Range("H25").Value = 26
$ws.Range("I25").Value = 0
$ws.Range("J25").Value = 1

$ws.Range("A26").Value = 'tostada con crema de cacao '
$ws.Range("B26").Value = 153.36
$ws.Range("C26").Value = 5.96
$ws.Range("D26").Value = 1.864
$ws.Range("E26").Value = 10.12
$ws.Range("F26").Value = 9.800000000000001
$ws.Range("G26").Value = 1.64
$ws.Range("H26").Value = 26
$ws.Range("I26").Value = 0
$ws.Range("J26").Value = 4

$ws.Range("A27").Value = '2 tostadas con creama de cacao'
$ws.Range("B27").Value = 300
$ws.Range("C27").Value = 11.9
$ws.Range("D27").Value = 2.5
$ws.Range("E27").Value = 20.24
$ws.Range("F27").Value = 19.6
$ws.Range("G27").Value = 3.28
$ws.Range("H27").Value = 26
$ws.Range("I27").Value = 0
$ws.Range("J27").Value = 4

$ws.Range("A28").Value = '150 gr pure de patata + 1 salchicha frankfurt + yogurt'
$ws.Range("B28").Value = 1092.84
$ws.Range("C28").Value = 32.08199999999999
$ws.Range("D28").Value = 12.147
$ws.Range("E28").Value = 174.41
$ws.Range("F28").Value = 51.129
$ws.Range("G28").Value = 43.664
$ws.Range("H28").Value = 5
$ws.Range("I28").Value = 0
$ws.Range("J28").Value = 3

$ws.Range("A29").Value = '150 gramos pure de patata + 100 gramos pollo frito + manzana'
$ws.Range("B29").Value = 883.5
$ws.Range("C29").Value = 21.76
$ws.Range("D29").Value = 2.85
$ws.Range("E29").Value = 142.94
$ws.Range("F29").Value = 0.2
$ws.Range("G29").Value = 24.88
$ws.Range("H29").Value = 5
$ws.Range("I29").Value = 0
$ws.Range("J29").Value = 1
